$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 142.5
$ws.Range("I33").Value = 180.2
$ws.Range("J33").Value = 115.57143
$ws.Range("K33").Value = 180.2
$ws.Range("L33").Value = 115.57143
$ws.Range("M33").Value = 48.80000000000001
$ws.Range("N33").Value = -573.57143

$ws.Range("H116").Value = 14584.167
$ws.Range("J116").Value = 7890
$ws.Range("L116").Value = 7890
$ws.Range("N116").Value = -14774

$ws.Range("H129").Value = 856.38776
$ws.Range("J129").Value = 925.14636
$ws.Range("L129").Value = 2775.43908
$ws.Range("N129").Value = -12775.43908

$ws.Range("H132").Value = 5465622
$ws.Range("I132").Value = 6411355.5
$ws.Range("J132").Value = 1381.7778
$ws.Range("K132").Value = 19234066.5
$ws.Range("L132").Value = 4145.3334
$ws.Range("M132").Value = -19231536.5
$ws.Range("N132").Value = -9205.3334

$ws.Range("H135").Value = 437.20587
$ws.Range("I135").Value = 398.90625
$ws.Range("J135").Value = 1050
$ws.Range("K135").Value = 3590.15625
$ws.Range("L135").Value = 9450
$ws.Range("M135").Value = -1055.15625
$ws.Range("N135").Value = -14520

$ws.Range("H141").Value = 637469.9
$ws.Range("I141").Value = 700704.4
$ws.Range("K141").Value = 2102113.2
$ws.Range("M141").Value = -2096933.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3226.9639
$ws.Range("I32").Value = 2576.4246
$ws.Range("K32").Value = 2576.4246
$ws.Range("M32").Value = -2289.4246

$ws.Range("H45").Value = 1743.7222
$ws.Range("I45").Value = 1574.375
$ws.Range("K45").Value = 1574.375
$ws.Range("M45").Value = -1197.375

$ws.Range("H61").Value = 40000828
$ws.Range("I61").Value = 20834196
$ws.Range("K61").Value = 20834196
$ws.Range("M61").Value = -20833984

$ws.Range("H74").Value = 1182.925
$ws.Range("I74").Value = 987.7037
$ws.Range("K74").Value = 987.7037
$ws.Range("M74").Value = -113.7037

$ws.Range("H77").Value = 1182.925
$ws.Range("I77").Value = 987.7037
$ws.Range("K77").Value = 4938.5185
$ws.Range("M77").Value = -570.5185000000001

$ws.Range("H136").Value = 40000828
$ws.Range("I136").Value = 20834196
$ws.Range("K136").Value = 62502588
$ws.Range("M136").Value = -62500038

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 304.2857
$ws.Range("I22").Value = 246.66667
$ws.Range("J22").Value = 347.5
$ws.Range("K22").Value = 246.66667
$ws.Range("L22").Value = 347.5
$ws.Range("M22").Value = -73.66667000000001
$ws.Range("N22").Value = -693.5

$ws.Range("H134").Value = 3670.3193
$ws.Range("I134").Value = 3056.0327
$ws.Range("K134").Value = 9168.098100000001
$ws.Range("M134").Value = -6633.098100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1544.4445
$ws.Range("I22").Value = 633.3333
$ws.Range("K22").Value = 633.3333
$ws.Range("M22").Value = -283.3333

$ws.Range("H31").Value = 1985575
$ws.Range("I31").Value = 3247706.8
$ws.Range("K31").Value = 3247706.8
$ws.Range("M31").Value = -3247411.8

$ws.Range("H34").Value = 1985575
$ws.Range("I34").Value = 3247706.8
$ws.Range("K34").Value = 3247706.8
$ws.Range("M34").Value = -3247504.8

$ws.Range("H58").Value = 714244.9399999999
$ws.Range("I58").Value = 1012058.2
$ws.Range("K58").Value = 1012058.2
$ws.Range("M58").Value = -1011855.2

$ws.Range("H99").Value = 2602.875
$ws.Range("I99").Value = 2080.75
$ws.Range("K99").Value = 2080.75
$ws.Range("M99").Value = -582.75

$ws.Range("H107").Value = 2337.3333
$ws.Range("I107").Value = 5000
$ws.Range("K107").Value = 5000
$ws.Range("M107").Value = -3080

$ws.Range("H126").Value = 2602.875
$ws.Range("I126").Value = 2080.75
$ws.Range("K126").Value = 6242.25
$ws.Range("M126").Value = -3772.25

$ws.Range("H132").Value = 1341.72
$ws.Range("I132").Value = 776.8
$ws.Range("J132").Value = 3601.4
$ws.Range("K132").Value = 2330.4
$ws.Range("L132").Value = 10804.2
$ws.Range("M132").Value = 199.6000000000004
$ws.Range("N132").Value = -15864.2

$ws.Range("H134").Value = 1327.3889
$ws.Range("I134").Value = 1077.0222
$ws.Range("K134").Value = 3231.0666
$ws.Range("M134").Value = -696.0666000000001

$ws.Range("H136").Value = 714244.9399999999
$ws.Range("I136").Value = 1012058.2
$ws.Range("K136").Value = 3036174.6
$ws.Range("M136").Value = -3033624.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 803
$ws.Range("I5").Value = 618.75
$ws.Range("J5").Value = 876.7
$ws.Range("K5").Value = 1856.25
$ws.Range("L5").Value = 2630.1
$ws.Range("M5").Value = -1744.25
$ws.Range("N5").Value = -2854.1

$ws.Range("H132").Value = 1118.8889
$ws.Range("I132").Value = 859.4
$ws.Range("J132").Value = 1443.25
$ws.Range("K132").Value = 7734.599999999999
$ws.Range("L132").Value = 12989.25
$ws.Range("M132").Value = -5204.599999999999
$ws.Range("N132").Value = -18049.25

$ws.Range("H135").Value = 803
$ws.Range("I135").Value = 618.75
$ws.Range("J135").Value = 876.7
$ws.Range("K135").Value = 5568.75
$ws.Range("L135").Value = 7890.3
$ws.Range("M135").Value = -3033.75
$ws.Range("N135").Value = -12960.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 676655.4
$ws.Range("I132").Value = 895728.75
$ws.Range("K132").Value = 2687186.25
$ws.Range("M132").Value = -2684656.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1809.2693
$ws.Range("I61").Value = 1565.421
$ws.Range("J61").Value = 2471.1428
$ws.Range("K61").Value = 1565.421
$ws.Range("L61").Value = 2471.1428
$ws.Range("M61").Value = -1363.421
$ws.Range("N61").Value = -2875.1428

$ws.Range("H68").Value = 2205.25
$ws.Range("I68").Value = 1773.6666
$ws.Range("K68").Value = 1773.6666
$ws.Range("M68").Value = -1024.6666

$ws.Range("H71").Value = 2205.25
$ws.Range("I71").Value = 1773.6666
$ws.Range("K71").Value = 8868.333000000001
$ws.Range("M71").Value = -5124.333000000001

$ws.Range("H113").Value = 1809.2693
$ws.Range("I113").Value = 1565.421
$ws.Range("J113").Value = 2471.1428
$ws.Range("K113").Value = 1565.421
$ws.Range("L113").Value = 2471.1428
$ws.Range("M113").Value = 604.579
$ws.Range("N113").Value = -6811.1428

$ws.Range("H132").Value = 1224.55
$ws.Range("I132").Value = 850
$ws.Range("K132").Value = 2550
$ws.Range("M132").Value = -20

$ws.Range("H136").Value = 1594.3529
$ws.Range("I136").Value = 886.2909
$ws.Range("K136").Value = 2658.8727
$ws.Range("M136").Value = -108.8726999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 28450.896
$ws.Range("I122").Value = 30417.629
$ws.Range("K122").Value = 91252.887
$ws.Range("M122").Value = -88802.887

$ws.Range("H132").Value = 1128.0366
$ws.Range("I132").Value = 718.4516
$ws.Range("J132").Value = 2397.75
$ws.Range("K132").Value = 2155.3548
$ws.Range("L132").Value = 7193.25
$ws.Range("M132").Value = 374.6451999999999
$ws.Range("N132").Value = -12253.25

$ws.Range("H136").Value = 8418838
$ws.Range("I136").Value = 10894497
$ws.Range("J136").Value = 1594.6666
$ws.Range("K136").Value = 32683491
$ws.Range("M136").Value = -32680941
$ws.Range("N136").Value = -9883.9998
